$d = $word.ActiveDocument

# "Pendiente, hasta el domingo a las 00:00" -> "HECHO"
$d.Content.Find.Execute("Pendiente, hasta el domingo a las 00:00", $true, $false, $false, $false, $false, $true, 1, $false, "HECHO", 2)

# "Pendiente, lunes antes de las 6 pm" -> "HECHO" (appears twice; wdReplaceAll handles both)
$d.Content.Find.Execute("Pendiente, lunes antes de las 6 pm", $true, $false, $false, $false, $false, $true, 1, $false, "HECHO", 2)

# "Pendiente, lunes" -> "HECHO"
$d.Content.Find.Execute("Pendiente, lunes", $true, $false, $false, $false, $false, $true, 1, $false, "HECHO", 2)
